# Moving from 2.0.2 to 2.0.3.
#
# The source change only touches the internal identifiers that Word mints
# for the "bookmark1" bookmark (w:bookmarkStart/@w:id, w:bookmarkEnd/@w:id)
# and the internal rsid stamp shared by the runs that make up the
# "REF bookmark1 \h" field right after it (w:r/@w:rsidR). These are
# regenerated housekeeping values produced whenever the fixture is
# rebuilt - no visible text, formatting, or document structure changes.
#
# The Word object model does not expose a way to poke an arbitrary literal
# value into those internal id/rsid slots (they are owned/assigned by
# Word itself), so the faithful way to reproduce "this bookmark's
# identity was refreshed" through COM automation is to recreate the
# bookmark in place: delete it and add it back around the exact same
# range/name. That forces Word to mint a fresh internal bookmark id for
# "bookmark1", mirroring the regeneration captured by the diff, while
# leaving the bookmarked text and every other part of the document
# untouched.

$d = $word.ActiveDocument

$bookmarkName = "bookmark1"
$bm = $d.Bookmarks($bookmarkName)

# Capture the exact extent of the existing bookmark before removing it.
$bmRange = $d.Range($bm.Start, $bm.End)

# Drop the old bookmark (carrying the stale identifier) ...
$bm.Delete()

# ... and recreate it on the same text span, so Word assigns it a brand
# new internal id while the bookmarked content stays identical.
$d.Bookmarks.Add($bookmarkName, $bmRange)
